# edit.ps1
# Applies the "Omit all GIS extensions except .shp and .dbf" change plus
# related proofing-mark cleanups described by the commit diff.

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p.Range
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: "Identical File" table cell description - split the single
# run into several runs (adding the GIS-extensions sentence, with
# en-CA language + spell-check markers on "shp"/"dbf"), and split the
# trailing example sentence into its own new paragraph labelled
# "Example 1:".
# ---------------------------------------------------------------------
$target = Get-ParagraphByText "*Flags duplicate files. Error count is incremented for each duplicate found. For example, if a group of 5 identical files are found, the error count is incremented by 4. (Includes owner column.)*"
if ($target -ne $null) {
    $target.Text = ""
    $xml = @"
<w:p $W>
  <w:r><w:t>Flags duplicate files. Error count is incremented for each duplicate found.</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>Most GIS extensions are omitted from output, except for .</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>shp</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> and .</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>dbf</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>.</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>(Includes owner column.)</w:t></w:r>
</w:p>
<w:p $W>
  <w:r><w:t>Example 1:</w:t></w:r>
  <w:r><w:t xml:space="preserve"> if a group of 5 identical files are found, the error count is incremented by 4.</w:t></w:r>
</w:p>
"@
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Change 2: "Bad Character (DIR)" table cell description - drop the
# gramStart/gramEnd proofing marks around "nor", merge it back into the
# surrounding sentence, and add a lastRenderedPageBreak marker.
# ---------------------------------------------------------------------
$target = Get-ParagraphByText "*Flags directory names with bad characters. A bad character is any character that is either not alphanumeric, nor a hyphen (-).*"
if ($target -ne $null) {
    $target.Text = ""
    $xml = @"
<w:p $W>
  <w:r><w:t xml:space="preserve">Flags directory names with bad characters. A bad character is </w:t></w:r>
  <w:r><w:lastRenderedPageBreak/><w:t>any character that is either not alphanumeric, nor a hyphen (-).</w:t></w:r>
</w:p>
"@
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Change 3: "Middle-click" hint bullet - drop gramStart/gramEnd proofing
# marks and merge the two runs into one.
# ---------------------------------------------------------------------
$target = Get-ParagraphByText "*Middle-click anywhere in the window to alternate between light and dark mode.*"
if ($target -ne $null) {
    $target.Text = ""
    $xml = @"
<w:p $W>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="5"/>
    </w:numPr>
  </w:pPr>
  <w:r><w:t>Middle-click anywhere in the window to alternate between light and dark mode.</w:t></w:r>
</w:p>
"@
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Change 4: "ctrl+w" hint bullet - add spellStart/spellEnd proofing marks
# around the verbatim "ctrl+w" run (keeping its VerbatimChar character
# style and the trailing bookmarkEnd).
# ---------------------------------------------------------------------
$target = Get-ParagraphByText "*ctrl+w to close the window.*"
if ($target -ne $null) {
    $target.Text = ""
    $xml = @"
<w:p $W>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="5"/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>ctrl+w</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> to close the window.</w:t></w:r>
  <w:bookmarkEnd w:id="5"/>
</w:p>
"@
    $target.InsertXML($xml)

    $styleRange = $d.Content
    $styleRange.Find.ClearFormatting()
    $styleRange.Find.Execute("ctrl+w")
    $styleRange.Style = "Verbatim Char"
}
